$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1, matching the style of the other header cells (H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Fill in the new I/J column values for each data row
$data = @(
    @(7, 8),
    @(8, 9),
    @(7, 7),
    @(8, 8),
    @(6, 7),
    @(1, 2),
    @(7, 8),
    @(4, 5),
    @(1, 2),
    @(2, 2),
    @(1, 1)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $pair = $data[$i]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
